# Updates cryptos list: refresh Price (col D) and Volume(1h) (col E) values
# for the rows whose figures moved, plus the Coin/Link/Price/Volume cells for
# rows 44-45 whose ranking order changed (dogwifhat <-> USDe).
# Leading "'" forces Excel to keep these as literal text (inlineStr) instead
# of auto-coercing number-looking strings like "1.00" into numeric values,
# matching the workbook's original all-text column layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.573.63"
$ws.Range("E2").Value = "'  -2.22%  "
$ws.Range("D3").Value = "'2.672.64"
$ws.Range("E3").Value = "'  -2.83%  "
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'600.46"
$ws.Range("E5").Value = "'  -0.69%  "
$ws.Range("D6").Value = "'166.89"
$ws.Range("E6").Value = "'  -0.08%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.545"
$ws.Range("E8").Value = "'  -0.30%  "
$ws.Range("D9").Value = "'2.672.02"
$ws.Range("E9").Value = "'  -2.86%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = "'  +0.87%  "
$ws.Range("D11").Value = "'0.158"
$ws.Range("D12").Value = "'0.362"
$ws.Range("E12").Value = "'  -0.65%  "
$ws.Range("D13").Value = "'5.21"
$ws.Range("E13").Value = "'  -2.56%  "
$ws.Range("D14").Value = "'27.87"
$ws.Range("E14").Value = "'  -3.53%  "
$ws.Range("D15").Value = "'3.167.16"
$ws.Range("E15").Value = "'  -2.59%  "
$ws.Range("D16").Value = "'0.0000184"
$ws.Range("E16").Value = "'  -3.45%  "
$ws.Range("D17").Value = "'67.521.77"
$ws.Range("E17").Value = "'  -2.17%  "
$ws.Range("D18").Value = "'2.666.56"
$ws.Range("E18").Value = "'  -2.32%  "
$ws.Range("D19").Value = "'11.72"
$ws.Range("E19").Value = "'  -2.32%  "
$ws.Range("D20").Value = "'7.85"
$ws.Range("E20").Value = "'  +1.21%  "
$ws.Range("D21").Value = "'364.56"
$ws.Range("E21").Value = "'  -1.24%  "
$ws.Range("D22").Value = "'4.39"
$ws.Range("E22").Value = "'  -4.33%  "
$ws.Range("D23").Value = "'4.81"
$ws.Range("E23").Value = "'  -3.76%  "
$ws.Range("E24").Value = "'  -5.07%  "
$ws.Range("D26").Value = "'70.82"
$ws.Range("E26").Value = "'  -4.38%  "
$ws.Range("D27").Value = "'10.13"
$ws.Range("E27").Value = "'  +1.39%  "
$ws.Range("D28").Value = "'2.820.58"
$ws.Range("E28").Value = "'  -1.96%  "
$ws.Range("D29").Value = "'0.0000102"
$ws.Range("E29").Value = "'  -4.48%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  +0.07%  "
$ws.Range("D31").Value = "'551.88"
$ws.Range("E31").Value = "'  -8.26%  "
$ws.Range("D32").Value = "'8.01"
$ws.Range("E32").Value = "'  -4.24%  "
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = "'  -4.80%  "
$ws.Range("E34").Value = "'  -2.35%  "
$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "'  -2.19%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  -0.06%  "
$ws.Range("D37").Value = "'1.55"
$ws.Range("E37").Value = "'  -5.52%  "
$ws.Range("D38").Value = "'19.48"
$ws.Range("E38").Value = "'  -3.57%  "
$ws.Range("D39").Value = "'155.11"
$ws.Range("E39").Value = "'  -5.07%  "
$ws.Range("D40").Value = "'0.372"
$ws.Range("E40").Value = "'  -2.90%  "
$ws.Range("E41").Value = "'  -4.94%  "
$ws.Range("D42").Value = "'5.28"
$ws.Range("E42").Value = "'  -4.64%  "
$ws.Range("D43").Value = "'17.92"
$ws.Range("E43").Value = "'  -0.60%  "
$ws.Range("B44").Value = "'USDe"
$ws.Range("C44").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("B45").Value = "'dogwifhat"
$ws.Range("C45").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.51"
$ws.Range("E45").Value = "'  -7.84%  "
$ws.Range("D46").Value = "'40.34"
$ws.Range("E46").Value = "'  -1.21%  "
$ws.Range("E47").Value = "'  -6.39%  "
$ws.Range("D48").Value = "'0.589"
$ws.Range("E48").Value = "'  -3.75%  "
$ws.Range("D49").Value = "'153.72"
$ws.Range("E49").Value = "'  -3.52%  "
$ws.Range("D50").Value = "'3.87"
$ws.Range("E50").Value = "'  -2.79%  "
$ws.Range("D51").Value = "'1.72"
$ws.Range("E51").Value = "'  -4.52%  "
